# Update diagrams in DG
# LogicComponentClassDiagram.pptx - rename a couple of class-diagram labels
# on the single slide:
#   "AddressBook" -> "BookShelf"   (Rectangle 62, shape id 16)
#   "FindCommand" -> "ListCommand" (Folded Corner 126, shape id 127 -
#                                   inside "XYZCommand = AddCommand, FindCommand, etc.")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

function Replace-TextRangeSubstring($textRange, [string]$oldText, [string]$newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -ge 0) {
        $part = $textRange.Characters($idx + 1, $oldText.Length)
        $part.Text = $newText
    }
}

# "AddressBook" -> "BookShelf"
$addressBookShape = Get-ShapeById $s 16
Replace-TextRangeSubstring $addressBookShape.TextFrame.TextRange "AddressBook" "BookShelf"

# "FindCommand" -> "ListCommand"
$findCommandShape = Get-ShapeById $s 127
Replace-TextRangeSubstring $findCommandShape.TextFrame.TextRange "FindCommand" "ListCommand"
